# Helper: build the OLE "long" RGB value (0x00BBGGRR) the way VBA's RGB() does,
# since this host does not expose the RGB() builtin.
function RGBVal($r, $g, $b) { return $r + ($g * 256) + ($b * 65536) }

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# 1) Table on slide 5 ("B1- TYPES OF FINANCIAL DOCUMENTS ") switches table
#    style from the default "No Style, Table Grid" to "No Style, No Grid".
# ---------------------------------------------------------------------------
for ($si = 1; $si -le $p.Slides.Count; $si++) {
    $slide = $p.Slides.Item($si)
    for ($shi = 1; $shi -le $slide.Shapes.Count; $shi++) {
        $shape = $slide.Shapes.Item($shi)
        if ($shape.HasTable) {
            $shape.Table.ApplyStyle("{0A967A63-D986-47C7-9DCC-72701F05C5F4}")
        }
    }
}

# ---------------------------------------------------------------------------
# 2) Swap the deck's colour theme from "Integral" (Red Violet) to the
#    stock "Office Theme" (Office) palette. dk1/lt1 (black/white) are
#    unchanged; the remaining ten theme colours are updated to match.
# ---------------------------------------------------------------------------
$theme = $p.SlideMaster.Theme
$colorScheme = $theme.ThemeColorScheme

# 1=dk1 2=lt1 3=dk2 4=lt2 5=accent1 6=accent2 7=accent3 8=accent4
# 9=accent5 10=accent6 11=hlink 12=folHlink  (standard MsoThemeColorSchemeIndex order)
$colorScheme.Item(3).RGB  = RGBVal 0x44 0x54 0x6A   # dk2
$colorScheme.Item(4).RGB  = RGBVal 0xE7 0xE6 0xE6   # lt2
$colorScheme.Item(5).RGB  = RGBVal 0x5B 0x9B 0xD5   # accent1
$colorScheme.Item(6).RGB  = RGBVal 0xED 0x7D 0x31   # accent2
$colorScheme.Item(7).RGB  = RGBVal 0xA5 0xA5 0xA5   # accent3
$colorScheme.Item(8).RGB  = RGBVal 0xFF 0xC0 0x00   # accent4
$colorScheme.Item(9).RGB  = RGBVal 0x44 0x72 0xC4   # accent5
$colorScheme.Item(10).RGB = RGBVal 0x70 0xAD 0x47   # accent6
$colorScheme.Item(11).RGB = RGBVal 0x05 0x63 0xC1   # hlink
$colorScheme.Item(12).RGB = RGBVal 0x95 0x4F 0x72   # folHlink
